# Insert one new data row before the current row 406, shifting existing
# rows 406:424 down to 407:425 (matches dimension change A1:R424 -> A1:R425).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(406).Insert()

$row = 406
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 45041
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino ensalada"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 400
$ws.Cells.Item($row, 11).Value = 15000
$ws.Cells.Item($row, 12).Value = 16000
$ws.Cells.Item($row, 13).Value = 15500
$ws.Cells.Item($row, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 258
$ws.Cells.Item($row, 17).Value = 60
$ws.Cells.Item($row, 18).Value = "Hortaliza"
